# Update "想去人数" (interest count, column F) figures for two worksheets
# ("展览" and "全部类型") to reflect the latest scrape snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 79
$ws1.Range("F3").Value = 3947
$ws1.Range("F4").Value = 2323
$ws1.Range("F5").Value = 461
$ws1.Range("F10").Value = 25
$ws1.Range("F11").Value = 119
$ws1.Range("F12").Value = 1464
$ws1.Range("F14").Value = 2681
$ws1.Range("F15").Value = 183

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 79
$ws4.Range("F3").Value = 3947
$ws4.Range("F4").Value = 2323
$ws4.Range("F5").Value = 461
$ws4.Range("F11").Value = 25
$ws4.Range("F12").Value = 119
$ws4.Range("F15").Value = 1464
$ws4.Range("F17").Value = 2681
$ws4.Range("F18").Value = 183
